# Implantação da busca pela home com falha com try/catch
#
# - Renomeia a aba "buscarHomeFail" para "buscarHomeFalha"
# - Preenche a aba buscarHomeFalha com o cabeçalho "Modelo" e o valor
#   "Dell Vostro" (novo cenário de busca com falha), espelhando o layout
#   já usado em buscarHomeSucesso
# - Ajusta a seleção/largura de coluna/configuração de página da nova aba

$wb = $excel.ActiveWorkbook

$wsFalha    = $wb.Worksheets.Item(4)   # buscarHomeFail -> buscarHomeFalha
$wsSucesso  = $wb.Worksheets.Item(3)   # buscarHomeSucesso (modelo de referência)

# 1) Renomeia a planilha
$wsFalha.Name = "buscarHomeFalha"

# 2) Conteúdo: cabeçalho "Modelo" e o novo modelo "Dell Vostro"
$wsFalha.Range("A1").Value = "Modelo"
$wsFalha.Range("A2").Value = "Dell Vostro"

# 3) Formata A1 como cabeçalho: negrito, centralizado e com quebra de texto
$wsFalha.Range("A1").HorizontalAlignment = -4108   # xlCenter
$wsFalha.Range("A1").WrapText = $true
$wsFalha.Range("A1").Font.Bold = $true

# 4) Formata A2 reaproveitando o mesmo estilo usado em buscarHomeSucesso!A2
#    (alinhamento vertical centralizado + quebra de texto)
$wsSucesso.Range("A2").Copy()
$wsFalha.Range("A2").PasteSpecial(-4122)           # xlPasteFormats
$excel.CutCopyMode = $false

# 5) Largura da coluna A para acomodar o texto do modelo
$wsFalha.Columns.Item(1).ColumnWidth = 33.6

# 6) Configuração de página (igual às demais planilhas)
$wsFalha.PageSetup.PaperSize = 9        # xlPaperA4... (A4/Letter conforme config local)
$wsFalha.PageSetup.Orientation = 1      # xlPortrait

# 7) Seleções finais nas planilhas afetadas
$wsFalha.Range("A2").Select()
$wsSucesso.Range("A1:A2").Select()

# Reativa a planilha buscarHomeFalha (mantém a aba em foco)
$wsFalha.Activate()
